$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: Assigned User header moves up from string index 397 to 394 (value stays "Assigned User")
# Nothing to change here - header text is unchanged, just the underlying shared string index shifts,
# which the engine will recompute automatically.

# Row 2: remove Jurisdiction Path value ("USA, State 1"), Assigned User becomes numeric 378
$ws.Range("CR2").ClearContents()
$ws.Range("CS2").Value = 378

# Row 3: remove Assigned User numeric value (4)
$ws.Range("CS3").ClearContents()

# Row 4: add Assigned User numeric value (2)
$ws.Range("CS4").Value = 2

# Row 5: Assigned User numeric value changes from 4 to 83
$ws.Range("CS5").Value = 83

# Row 6: remove Jurisdiction Path value ("USA, State 1, County 1"), add Assigned User numeric value 83
$ws.Range("CR6").ClearContents()
$ws.Range("CS6").Value = 83

# Row 7: remove Jurisdiction Path value ("USA, State 1, County 2")
$ws.Range("CR7").ClearContents()

# Row 8: remove Jurisdiction Path value ("USA, State 1, County 2"), Assigned User changes from 2 to 57
$ws.Range("CR8").ClearContents()
$ws.Range("CS8").Value = 57

# Row 9: remove Assigned User numeric value (2)
$ws.Range("CS9").ClearContents()

# Row 10: add Assigned User numeric value (9999)
$ws.Range("CS10").Value = 9999

# Update the active selection on the sheet
$ws.Range("CS7").Select()
